# feat(excel): Update Read function
# Adds two new columns to the "江苏城市分级" sheet:
#   N: "是的" - a yes/no (boolean) flag per city row
#   O: "日期" - a date stamp per city row, formatted as yyyy/m/d;@

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$ws.Range("N1").Value = "是的"
$ws.Range("O1").Value = "日期"

# New boolean flag column (rows 2-14)
$flags = @($true, $true, $true, $false, $true, $true, $false, $false, $false, $true, $true, $true, $true)
$row = 2
foreach ($flag in $flags) {
    $ws.Cells.Item($row, 14).Value = $flag
    $row++
}

# New date column (rows 2-14), formatted as yyyy/m/d;@
$ws.Range("O2:O14").NumberFormat = "yyyy/m/d;@"
$year = 2024
$month = 1
$day = 9
for ($r = 2; $r -le 14; $r++) {
    $d = Get-Date -Year $year -Month $month -Day $day -Hour 0 -Minute 0 -Second 0
    $ws.Cells.Item($r, 15).Value = $d
    $day++
}

# Restore the last active selection
$ws.Range("Q18").Select() | Out-Null
